$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K4").Value = "Chardonnay, 750 ML"

$ws.Range("K5").Value = "Blueberry Wensleydale"
$ws.Range("L5").Value = 40.81632653061224

$ws.Range("K6").Value = "Cremant Rose, 750 ML"
$ws.Range("L6").Value = 43.13725490196079

$ws.Range("K7").Value = "Imperial Rioja Reserva, 750 ML"
$ws.Range("L7").Value = 42.30769230769231

$ws.Range("K10").Value = "Red Wine Vinegar, 16.9 FZ"
$ws.Range("L10").Value = 42.10526315789474

$ws.Range("K11").Value = "Vinegar Rice Premium, 10 FZ"
$ws.Range("L11").Value = 40

$ws.Range("K12").Value = "Vinegar Rice Premium, 10 FZ"
$ws.Range("L12").Value = 40

$ws.Range("K16").Value = "Apple Juice, 64 FZ"
$ws.Range("L16").Value = 36.36363636363637

$ws.Range("K18").Value = "Pecorino Romano"
$ws.Range("L18").Value = 34.48275862068965

$ws.Range("K19").Value = "Pecorino Romano"
$ws.Range("L19").Value = 34.48275862068965

$ws.Range("K20").Value = "Cinnamon Toast Cereal, 9 OZ"
$ws.Range("L20").Value = 48.38709677419355

$ws.Range("K21").Value = "Brut, 750 ML"
$ws.Range("L21").Value = 35.55555555555556

$ws.Range("K24").Value = "Taleggio"
$ws.Range("L24").Value = 40

$ws.Range("K25").Value = "Taleggio"
$ws.Range("L25").Value = 40
